# Append five new species-observation rows (7-11) to the Artfynd sheet,
# matching the columns used by the existing rows 2-6.
#
# Note: a handful of text values look like a number ("1") or a date
# ("2025-10-10"), and Excel would silently reinterpret those as numeric/
# date values if assigned via .Value. Writing them through .Formula with
# a leading "'" (the same trick Excel itself uses for "typed as text")
# keeps them as literal text, which is what the source data has. The
# same trick is used for cells that must exist but hold an empty string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 131321686
$ws.Range("B7").Value = 91776
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 5447
$ws.Range("F7").Value = 'Vedticka'
$ws.Range("G7").Value = 'Fuscoporia viticola'
$ws.Range("H7").Value = '(Schwein.) Murrill'
$ws.Range("I7").Formula = '''1'
$ws.Range("P7").Value = 'Bergsboda, Vb'
$ws.Range("Q7").Value = 762874
$ws.Range("R7").Value = 7082171
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = 'Västerbotten'
$ws.Range("U7").Value = 'Umeå'
$ws.Range("V7").Value = 'Västerbotten'
$ws.Range("W7").Value = 'Umeå stad'
$ws.Range("Y7").Formula = '''2025-10-10'
$ws.Range("Z7").Value = '11:13'
$ws.Range("AA7").Formula = '''2025-10-10'
$ws.Range("AB7").Value = '11:13'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Formula = ''''
$ws.Range("AW7").Value = 'Edvin Strandberg'
$ws.Range("AX7").Value = 'Edvin Strandberg'
$ws.Range("AY7").Formula = ''''

# Row 8
$ws.Range("A8").Value = 131321684
$ws.Range("B8").Value = 57884
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = 'Tretåig hackspett'
$ws.Range("G8").Value = 'Picoides tridactylus'
$ws.Range("H8").Value = '(Linnaeus, 1758)'
$ws.Range("I8").Formula = '''1'
$ws.Range("M8").Value = 'färska spår'
$ws.Range("P8").Value = 'Bergsboda, Vb'
$ws.Range("Q8").Value = 762920
$ws.Range("R8").Value = 7082182
$ws.Range("S8").Value = 5
$ws.Range("T8").Value = 'Västerbotten'
$ws.Range("U8").Value = 'Umeå'
$ws.Range("V8").Value = 'Västerbotten'
$ws.Range("W8").Value = 'Umeå stad'
$ws.Range("Y8").Formula = '''2025-10-10'
$ws.Range("Z8").Value = '10:59'
$ws.Range("AA8").Formula = '''2025-10-10'
$ws.Range("AB8").Value = '10:59'
$ws.Range("AC8").Value = 'Två granar, hyggeskanten.'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AT8").Formula = ''''
$ws.Range("AW8").Value = 'Edvin Strandberg'
$ws.Range("AX8").Value = 'Edvin Strandberg'
$ws.Range("AY8").Formula = ''''

# Row 9
$ws.Range("A9").Value = 131321677
$ws.Range("B9").Value = 91776
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 5447
$ws.Range("F9").Value = 'Vedticka'
$ws.Range("G9").Value = 'Fuscoporia viticola'
$ws.Range("H9").Value = '(Schwein.) Murrill'
$ws.Range("I9").Formula = '''1'
$ws.Range("P9").Value = 'Bergsboda, Vb'
$ws.Range("Q9").Value = 763042
$ws.Range("R9").Value = 7082241
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 'Västerbotten'
$ws.Range("U9").Value = 'Umeå'
$ws.Range("V9").Value = 'Västerbotten'
$ws.Range("W9").Value = 'Umeå stad'
$ws.Range("Y9").Formula = '''2025-10-10'
$ws.Range("Z9").Value = '10:31'
$ws.Range("AA9").Formula = '''2025-10-10'
$ws.Range("AB9").Value = '10:31'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AT9").Formula = ''''
$ws.Range("AW9").Value = 'Edvin Strandberg'
$ws.Range("AX9").Value = 'Edvin Strandberg'
$ws.Range("AY9").Formula = ''''

# Row 10
$ws.Range("A10").Value = 131321669
$ws.Range("B10").Value = 91776
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 5447
$ws.Range("F10").Value = 'Vedticka'
$ws.Range("G10").Value = 'Fuscoporia viticola'
$ws.Range("H10").Value = '(Schwein.) Murrill'
$ws.Range("I10").Formula = '''1'
$ws.Range("P10").Value = 'Bergsboda, Vb'
$ws.Range("Q10").Value = 763007
$ws.Range("R10").Value = 7082257
$ws.Range("S10").Value = 5
$ws.Range("T10").Value = 'Västerbotten'
$ws.Range("U10").Value = 'Umeå'
$ws.Range("V10").Value = 'Västerbotten'
$ws.Range("W10").Value = 'Umeå stad'
$ws.Range("Y10").Formula = '''2025-10-10'
$ws.Range("Z10").Value = '09:58'
$ws.Range("AA10").Formula = '''2025-10-10'
$ws.Range("AB10").Value = '09:58'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AT10").Formula = ''''
$ws.Range("AW10").Value = 'Edvin Strandberg'
$ws.Range("AX10").Value = 'Edvin Strandberg'
$ws.Range("AY10").Formula = ''''

# Row 11
$ws.Range("A11").Value = 131321688
$ws.Range("B11").Value = 91833
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = 'Granticka'
$ws.Range("G11").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H11").Formula = ''''
$ws.Range("I11").Formula = '''1'
$ws.Range("P11").Value = 'Bergsboda, Vb'
$ws.Range("Q11").Value = 762922
$ws.Range("R11").Value = 7082225
$ws.Range("S11").Value = 5
$ws.Range("T11").Value = 'Västerbotten'
$ws.Range("U11").Value = 'Umeå'
$ws.Range("V11").Value = 'Västerbotten'
$ws.Range("W11").Value = 'Umeå stad'
$ws.Range("Y11").Formula = '''2025-10-10'
$ws.Range("Z11").Value = '11:45'
$ws.Range("AA11").Formula = '''2025-10-10'
$ws.Range("AB11").Value = '11:45'
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AT11").Formula = ''''
$ws.Range("AW11").Value = 'Edvin Strandberg'
$ws.Range("AX11").Value = 'Edvin Strandberg'
$ws.Range("AY11").Formula = ''''
